$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Row, open_price(D), close_price(E), high_price(F), low_price(G), shares_outstanding(H)
$rowsData = @(
    ,@(2, 23.98999977111816, 36.5, 41.88999938964844, 23.65999984741211, 153433606)
    ,@(3, 23.98999977111816, 36.5, 41.88999938964844, 23.65999984741211, 153433606)
    ,@(4, 23.98999977111816, 36.5, 41.88999938964844, 23.65999984741211, 153433606)
    ,@(5, 23.98999977111816, 36.5, 41.88999938964844, 23.65999984741211, 153433606)
    ,@(6, 23.98999977111816, 36.5, 41.88999938964844, 23.65999984741211, 153433606)
    ,@(7, 36.97000122070312, 40.63000106811523, 44.79999923706055, 33.06999969482422, 153433606)
    ,@(8, 62.5, 34.11999893188477, 66.63999938964844, 33.75, 153433606)
    ,@(9, 29.25, 28.82999992370605, 30.5, 25.97999954223633, 153433606)
    ,@(10, 28.76000022888184, 33.04999923706055, 33.36000061035156, 26.76000022888184, 153433606)
    ,@(11, 29.10000038146973, 29.17000007629395, 32.22999954223633, 28.09000015258789, 153433606)
    ,@(12, 29.8700008392334, 31.95000076293945, 33.06900024414062, 29.85000038146973, 153433606)
    ,@(13, 24.1299991607666, 26.23999977111816, 27.45999908447266, 24.10000038146973, 153433606)
    ,@(14, 37.84999847412109, 42.20999908447266, 42.88000106811523, 36.06999969482422, 153433606)
    ,@(15, 55.25, 57.88999938964844, 65.37000274658203, 54.5009994506836, 153433606)
    ,@(16, 87, 75.22000122070312, 87.55000305175781, 62.43000030517578, 153433606)
    ,@(17, 87.44999694824219, 111.3199996948242, 113.3249969482422, 80.12999725341797, 153433606)
    ,@(18, 131.5800018310547, 137.1399993896484, 137.3999938964844, 118.4800033569336, 153433606)
    ,@(19, 138.3200073242188, 139.1100006103516, 150.4400024414062, 134.7449951171875, 153433606)
    ,@(20, 109.5999984741211, 96.55999755859376, 118.3399963378906, 92.88999938964844, 153433606)
    ,@(21, 100.0100021362305, 124.3399963378906, 127.5, 99.66000366210938, 153433606)
    ,@(22, 87, 112.3000030517578, 116.8399963378906, 79.25, 153433606)
    ,@(23, 218.1300048828125, 277.4200134277344, 277.4200134277344, 214.3200073242188, 153433606)
    ,@(24, 249.9499969482422, 278.9700012207031, 341.7000122070312, 247.5, 153433606)
    ,@(25, 339, 359.4299926757812, 404.2699890136719, 326.1099853515625, 153433606)
    ,@(26, 350.3599853515625, 367.7999877929688, 405.7999877929688, 336.3580017089844, 153433606)
    ,@(27, 396, 373.5899963378906, 412.6799926757813, 361.2799987792969, 153433606)
    ,@(28, 319.75, 291.3599853515625, 373, 278, 153433606)
    ,@(29, 264, 206.1199951171875, 266.739990234375, 172.6100006103516, 153433606)
    ,@(30, 166.0399932861328, 111.8199996948242, 176.3999938964844, 109.6790008544922, 153433606)
    ,@(31, 85.23000335693359, 84.80000305175781, 95.94000244140624, 78.41999816894531, 153433606)
    ,@(32, 69.51000213623047, 74.37000274658203, 78.43000030517578, 62.09999847412109, 153433606)
    ,@(33, 50, 59.84000015258789, 61.49200057983398, 48.04999923706055, 153433606)
    ,@(34, 65.34999847412109, 52.61000061035156, 66, 51.38999938964844, 153433606)
    ,@(35, 63.95999908447266, 66.02999877929688, 71.23999786376953, 59.87099838256836, 153433606)
    ,@(36, 58.33000183105469, 51.2599983215332, 58.83000183105469, 49.90999984741211, 153433606)
    ,@(37, 74.36000061035156, 70.33000183105469, 76.02999877929688, 67.73999786376953, 153433606)
    ,@(38, 61.34999847412109, 59.88000106811523, 63.40000152587891, 57.5, 153433606)
    ,@(39, 57, 59.13000106811523, 60.65999984741211, 55.38000106811523, 153433606)
    ,@(40, 65.19999694824219, 80.65000152587891, 83.30000305175781, 63.18999862670898, 153433606)
    ,@(41, 109.5500030517578, 146.5800018310547, 151.9499969482422, 105.5299987792969, 153433606)
    ,@(42, 98.08999633789062, 96.70999908447266, 101.4000015258789, 77.51200103759766, 153433606)
    ,@(43, 124.5599975585938, 129, 135.3999938964844, 112.1100006103516, 153433606)
)

foreach ($row in $rowsData) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $ws.Cells.Item($r, 8).Value = $row[5]
    $ws.Cells.Item($r, 9).Value = "TWLO"
}
